$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.439.94'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').Value = '3.249.38'
$ws.Range('E3').Value = '  +3.12%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'594.18"
$ws.Range('E5').Value = '  -1.54%  '
$ws.Range('D6').Value = "'141.83"
$ws.Range('E6').Value = '  -1.40%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.244.94'
$ws.Range('E8').Value = '  +3.16%  '
$ws.Range('E9').Value = '  -1.10%  '
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('D11').Value = "'5.36"
$ws.Range('E11').Value = '  -0.64%  '
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('E13').Value = '  -2.68%  '
$ws.Range('D14').Value = "'34.43"
$ws.Range('E14').Value = '  -1.66%  '
$ws.Range('D15').Value = '3.783.94'
$ws.Range('E15').Value = '  +3.10%  '
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').Value = '3.247.31'
$ws.Range('E17').Value = '  +2.78%  '
$ws.Range('D18').Value = '63.422.53'
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').Value = "'6.79"
$ws.Range('E19').Value = '  -1.06%  '
$ws.Range('D20').Value = "'477.29"
$ws.Range('E20').Value = '  -2.72%  '
$ws.Range('E22').Value = '  +1.61%  '
$ws.Range('D23').Value = "'7.95"
$ws.Range('E23').Value = '  +3.64%  '
$ws.Range('D24').Value = "'84.14"
$ws.Range('E24').Value = '  -4.34%  '
$ws.Range('D25').Value = "'13.22"
$ws.Range('E25').Value = '  -0.65%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = "'7.46"
$ws.Range('E27').Value = '  +6.72%  '
$ws.Range('D29').Value = "'8.07"
$ws.Range('E29').Value = '  -1.90%  '
$ws.Range('E30').Value = '  +3.07%  '
$ws.Range('D31').Value = "'27.64"
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  -2.51%  '
$ws.Range('D34').Value = "'2.55"
$ws.Range('E34').Value = '  -4.08%  '
$ws.Range('E35').Value = '  -1.61%  '
$ws.Range('D36').Value = "'5.91"
$ws.Range('E36').Value = '  -2.21%  '
$ws.Range('D37').Value = "'52.91"
$ws.Range('E37').Value = '  +0.41%  '
$ws.Range('E38').Value = '  -3.74%  '
$ws.Range('E39').Value = '  -0.81%  '
$ws.Range('D40').Value = "'420.12"
$ws.Range('E40').Value = '  -3.24%  '
$ws.Range('D41').Value = "'8.40"
$ws.Range('E41').Value = '  +1.00%  '
$ws.Range('D42').Value = '2.979.48'
$ws.Range('E42').Value = '  +1.48%  '
$ws.Range('D43').Value = "'2.76"
$ws.Range('E43').Value = '  -7.16%  '
$ws.Range('D44').Value = "'0.111"
$ws.Range('E44').Value = '  -7.51%  '
$ws.Range('E45').Value = '  +3.09%  '
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').Value = "'25.93"
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('D49').Value = "'2.33"
$ws.Range('E49').Value = '  -3.35%  '
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('D51').Value = "'120.82"
$ws.Range('E51').Value = '  +0.36%  '
